$d = $word.ActiveDocument

# Fix the bio paragraph for Steve Deakin-Davies
$d.Content.Find.Execute("having living nearby previously", $true, $false, $false, $false, $false, $true, 1, $false, "having lived nearby previously", 2)
$d.Content.Find.Execute("develop a whole range of businesses", $true, $false, $false, $false, $false, $true, 1, $false, "develop a range of businesses", 2)
$d.Content.Find.Execute("has a particular interest in health matters", $true, $false, $false, $false, $false, $true, 1, $false, "has an interest in health matters", 2)

# Move the _GoBack bookmark to mark the last edit location, as Word does automatically:
# position it right after "...chairs and v" in the Annual Meeting paragraph (mid-word,
# splitting "vice-chairs").
$rng = $d.Content
$rng.Find.Execute("chairs and v") | Out-Null
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null
